$d = $word.ActiveDocument

# --- Clean up spell-check (proofErr) run splits by re-asserting the full,
# --- already-correct sentence text. Word's Find/Replace collapses a
# --- paragraph's runs into a single run (dropping any <w:proofErr/> markers
# --- that were splitting it) when the replacement text spans the whole run
# --- set with matching content.

$d.Content.Find.Execute( `
    "Capabilities attributes such as hunger, what type of dinosaur can probably be implemented using enums in capabilities", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Capabilities attributes such as hunger, what type of dinosaur can probably be implemented using enums in capabilities", 2)

$d.Content.Find.Execute( `
    "Should accept probably player (enum maybe)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Should accept probably player (enum maybe)", 2)

$d.Content.Find.Execute( `
    "Will work on Hunger, Breeding and PlayerActions:", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Will work on Hunger, Breeding and PlayerActions:", 2)

$d.Content.Find.Execute( `
    "New JurassicWorld class to implement rain, general framework to be used by Amos", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "New JurassicWorld class to implement rain, general framework to be used by Amos", 2)

$d.Content.Find.Execute( `
    "Updating EatPreyAction and various things related to it", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Updating EatPreyAction and various things related to it", 2)

# --- Accept the WBA: add "I accept this WBA." right after "Lin Chen Xiang: "
# --- on the Assignment 3 sign-off block, as its own run (matching the
# --- sibling "Ng Yu Kang:" / "Amos Leong Zheng Khang:" acceptance lines).
# --- Inserting under Track Changes and then accepting the resulting
# --- revision keeps the new text in its own <w:r> instead of merging it
# --- into the existing "Lin Chen Xiang: " run.

$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true

foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text
    if ($ptext -eq "Lin Chen Xiang: `r") {
        $fr = $p.Range.Duplicate
        $fr.Find.Execute("Lin Chen Xiang: ")
        $fr.InsertAfter("I accept this WBA.")
    }
}

$d.TrackRevisions = $wasTracking
$d.Revisions.AcceptAll()
